$d = $word.ActiveDocument

# --- Hunk 1: "{{ p" + "lot" + "sHeaderLabel }}" -> single run "{{ plotsHeaderLabel }}" ---
# The text is unchanged once concatenated; Word's Find/Replace merges the
# three runs that carry it into a single run (taking the first run's
# formatting), which is exactly the target shape.
$target1 = "{{ plotsHeaderLabel }}"
$found1 = $d.Content.Find.Execute($target1, $true, $false, $false, $false, $false, $true, 1, $false, $target1, 2)
if (-not $found1) {
    throw "hunk1: target text not found"
}

# --- Hunk 2: "{% if p" + "lot" + "sHeader%}{{ p" + "lots" + "Header }}{% else %}-{% endif%}"
#             -> single run "{% if plotsHeader%}{{ plotsHeader }}{% else %}-{% endif%}" ---
# Same situation: concatenated text is unchanged, only the run split collapses.
$target2 = "{% if plotsHeader%}{{ plotsHeader }}{% else %}-{% endif%}"
$found2 = $d.Content.Find.Execute($target2, $true, $false, $false, $false, $false, $true, 1, $false, $target2, 2)
if (-not $found2) {
    throw "hunk2: target text not found"
}

# --- Hunk 3: the inputDateHeader cell -----------------------------------
# * <w:ilvl w:val="5"/> -> <w:ilvl w:val="2"/>
# * the single run's text is split into three runs, inserting a new
#   paperInputDateHeader conditional between the existing pieces.
# Locate the paragraph by its distinctive text and rewrite it in one shot
# via InsertXML (Find/Replace cannot touch pPr/numPr).
$target3 = "{% if inputDateHeader %}{{ inputDateHeader }}{% else %}-{% endif%}"
$found3 = $d.Content.Find.Execute($target3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "hunk3: target paragraph not found"
}

$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd([char]13, [char]7) -eq $target3) {
        $para = $candidate
        break
    }
}
if ($null -eq $para) {
    throw "hunk3: could not re-locate paragraph by text"
}

$newParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:suppressAutoHyphens w:val="true"/><w:spacing w:before="0" w:after="0"/><w:ind w:right="176" w:hanging="0"/><w:rPr><w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans" w:eastAsia="Noto Sans CJK SC Regular" w:cs="Arial"/><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="auto"/><w:kern w:val="2"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="fr-CH" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Arial" w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="auto"/><w:kern w:val="2"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="fr-CH" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>{% if inputDateHeader %}{{ inputDateHeader }}</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Arial" w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="auto"/><w:kern w:val="2"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="fr-CH" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>{% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Noto Sans CJK SC Regular" w:cs="Arial" w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="auto"/><w:kern w:val="2"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="fr-CH" w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>{% else %}-{% endif %}</w:t></w:r></w:p>
'@

$para.Range.InsertXML($newParaXml)
